# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (169-170) into the Mandarina sheet,
# pushing the existing rows 169-188 down to 171-190.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 169 - this shifts rows 169-188
# down to 171-190 and extends the used range/dimension automatically.
$ws.Rows("169:170").Insert()

# --- New row 169 ---------------------------------------------------------
$ws.Range("A169").Value = 11
$ws.Range("B169").Value = "Vega Monumental Concepción"
$ws.Range("C169").Value = "Bíobío"
$ws.Range("D169").Value = 45021
$ws.Range("E169").Value = 8
$ws.Range("F169").Value = "Fruta"
$ws.Range("G169").Value = 100102
$ws.Range("H169").Value = "Cítricos"
$ws.Range("I169").Value = 100102004
$ws.Range("J169").Value = "Mandarina"
$ws.Range("K169").Value = "Murcott"
$ws.Range("L169").Value = "Especial"
$ws.Range("M169").Value = 200
$ws.Range("N169").Value = 11000
$ws.Range("O169").Value = 12000
$ws.Range("P169").Value = 11500
$ws.Range("Q169").Value = "$/bandeja 18 kilos"
$ws.Range("R169").Value = "Provincia de Limarí"
$ws.Range("S169").Value = 639
$ws.Range("T169").Value = 18

# --- New row 170 ---------------------------------------------------------
$ws.Range("A170").Value = 11
$ws.Range("B170").Value = "Vega Monumental Concepción"
$ws.Range("C170").Value = "Bíobío"
$ws.Range("D170").Value = 45021
$ws.Range("E170").Value = 8
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100102
$ws.Range("H170").Value = "Cítricos"
$ws.Range("I170").Value = 100102004
$ws.Range("J170").Value = "Mandarina"
$ws.Range("K170").Value = "Murcott"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 250
$ws.Range("N170").Value = 10000
$ws.Range("O170").Value = 10000
$ws.Range("P170").Value = 10000
$ws.Range("Q170").Value = "$/bandeja 18 kilos"
$ws.Range("R170").Value = "Provincia de Limarí"
$ws.Range("S170").Value = 556
$ws.Range("T170").Value = 18
